$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by duplicating the existing
#    "2022-Q1" sheet (this carries over identical column headers, number
#    formats, borders and page setup) and dropping it in front of
#    "2022-Q1", so the final order is:
#      总计, 2022-Q4, 2022-Q1, 2021-Q4, 2021-Q3
# ---------------------------------------------------------------------------
$q1sheet = $wb.Worksheets.Item("2022-Q1")
$q1sheet.Copy($q1sheet)
$q4sheet = $wb.Worksheets.Item("2022-Q1 (2)")
$q4sheet.Name = "2022-Q4"

# The duplicated sheet only has 2 data rows (rows 2-3); extend it with two
# more rows (4-5), copying the formatting already used on row 3 so the new
# cells inherit the same (unstyled/default) look as the other data rows.
$q4sheet.Range("A3:H3").Copy()
$q4sheet.Range("A4:H5").PasteSpecial(-4122) # xlPasteFormats

# Fund-holding table for 2022-Q4 (header row 1 already matches; fill rows
# 2-5). Columns B-G are text values (fund codes, names and percentages
# stored as text, matching the sibling quarter sheets); columns A and H are
# plain numbers.
$q4data = @(
    @("011748", "华泰柏瑞景气成长混合A", "2.22", "79.84", "2.27", "0.0504", 10),
    @("005409", "华泰柏瑞战略新兴产业混合A", "1.34", "82.34", "2.33", "0.0312", 10),
    @("010032", "华泰柏瑞战略新兴产业混合C", "0.12", "82.34", "2.33", "0.0028", 10),
    @("011749", "华泰柏瑞景气成长混合C", "0.09", "79.84", "2.27", "0.0020", 10)
)

$q4sheet.Range("B2:G5").NumberFormat = "@"
for ($r = 0; $r -lt $q4data.Length; $r++) {
    $row = $r + 2
    $q4sheet.Cells.Item($row, 1).Value = $r
    $q4sheet.Cells.Item($row, 2).Value = $q4data[$r][0]
    $q4sheet.Cells.Item($row, 3).Value = $q4data[$r][1]
    $q4sheet.Cells.Item($row, 4).Value = $q4data[$r][2]
    $q4sheet.Cells.Item($row, 5).Value = $q4data[$r][3]
    $q4sheet.Cells.Item($row, 6).Value = $q4data[$r][4]
    $q4sheet.Cells.Item($row, 7).Value = $q4data[$r][5]
    $q4sheet.Cells.Item($row, 8).Value = $q4data[$r][6]
}

# The explicit "@" text number-format above leaves a style trace on the
# cells; wipe it by pasting in the plain/default formatting of an empty
# cell, restoring the same (unstyled) look the sibling sheets use for their
# data cells.
$q4sheet.Range("Z99").Copy()
$q4sheet.Range("B2:G5").PasteSpecial(-4122) # xlPasteFormats
$q4sheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add the 2022-Q4 row above the
#    existing quarters, shifting them down, and refresh the sequential index
#    stored in column A.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:D5").Clear()

$totalData = @(
    @(0, "2022-Q4", 4, 0.09),
    @(1, "2022-Q1", 2, 0.18),
    @(2, "2021-Q4", 12, 5.56),
    @(3, "2021-Q3", 1, 1.59)
)
for ($i = 0; $i -lt $totalData.Length; $i++) {
    $row = $i + 2
    $total.Cells.Item($row, 1).Value = $totalData[$i][0]
    $total.Cells.Item($row, 2).Value = $totalData[$i][1]
    $total.Cells.Item($row, 3).Value = $totalData[$i][2]
    $total.Cells.Item($row, 4).Value = $totalData[$i][3]
}

# Re-apply the bold/bordered look used on column A (matches header style).
$total.Range("B1").Copy()
$total.Range("A2:A5").PasteSpecial(-4122) # xlPasteFormats
$total.Range("A1").Select()
